# Updated cryptos list data: refresh Price (D) and Volume(1h) (E) columns,
# and swap the RenderToken/Hedera rows (35/36) per upstream source re-ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.754.77"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.819.35"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.70"
$ws.Range("E5").Value = "  +5.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.89"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +3.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.47"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.89"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.73"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "3.254.68"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "2.824.23"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.884"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "51.593.34"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  +9.05%  "
$ws.Range("E20").Value = "  -4.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "0.0₃0992"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.59"
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.66"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.71"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.02"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.54"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0451"
$ws.Range("E33").Value = "  +27.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  +4.44%  "
# Row 35: was Hedera -> now RenderToken (re-ranked above Hedera)
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  +4.66%  "

# Row 36: was RenderToken -> now Hedera
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0823"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.13"
$ws.Range("E40").Value = "  -5.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.83"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.22"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "2.079.76"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.65"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  +6.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.75"
$ws.Range("E51").Value = "  +0.50%  "
